$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "584.52")
# round-trip exactly instead of being auto-coerced into floating point numbers
# (which would lose trailing zeros / introduce binary rounding error).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.982.78"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.259.41"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "584.52"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "184.52"
$ws.Range("E6").Value = "  +4.27%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "0.134"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("D12").Value = "3.823.56"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D14").Value = "28.57"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "67.981.60"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "3.262.42"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "5.87"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "382.40"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").Value = "7.70"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "71.40"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "0.514"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "9.89"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "0.183"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "5.70"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  +7.93%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").Value = "162.60"
$ws.Range("E36").Value = "  -6.10%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "0.835"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").Value = "26.69"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "6.72"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("D41").Value = "4.61"
$ws.Range("E41").Value = "  +6.71%  "
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").Value = "0.0689"
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("D46").Value = "2.645.94"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").Value = "344.73"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "31.36"
$ws.Range("E51").Value = "  +3.30%  "
